# Remove the space on building no at url
#
# The "building_no" column header is renamed to "No", and the mazemap
# hyperlink display text for C2:C8 has its leading/trailing spaces around
# the building number removed (e.g. "sharepoi= 411 " -> "sharepoi=411").
# The underlying hyperlink relationships/targets themselves are left as-is;
# only the cell text (shared-string) content changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "No"

$ws.Range("C2").Value = "https://use.mazemap.com/?campusid=217&sharepoitype=identifier&sharepoi=411"
$ws.Range("C3").Value = "https://use.mazemap.com/?campusid=217&sharepoitype=identifier&sharepoi=416"
$ws.Range("C4").Value = "https://use.mazemap.com/?campusid=217&sharepoitype=identifier&sharepoi=417"
$ws.Range("C5").Value = "https://use.mazemap.com/?campusid=217&sharepoitype=identifier&sharepoi=418"
$ws.Range("C6").Value = "https://use.mazemap.com/?campusid=217&sharepoitype=identifier&sharepoi=421"
$ws.Range("C7").Value = "https://use.mazemap.com/?campusid=217&sharepoitype=identifier&sharepoi=434"
$ws.Range("C8").Value = "https://use.mazemap.com/?campusid=217&sharepoitype=identifier&sharepoi=437"
